{"js": "// Replace the date line and each two-digit multiplication problem's text\n// with the updated values, matching the authoritative diff 1:1.\nconst body = context.document.body;\nconst pairs = [\n  [\"2024-12-15 Sunday\", \"2024-12-16 Monday\"],\n  [\"56\u00d794=5264\", \"39\u00d786=3354\"],\n  [\"32\u00d762=1984\", \"77\u00d775=5775\"],\n  [\"46\u00d746=2116\", \"55\u00d720=1100\"],\n  [\"34\u00d730=1020\", \"12\u00d711=132\"],\n  [\"25\u00d744=1100\", \"61\u00d724=1464\"],\n  [\"31\u00d732=992\", \"31\u00d752=1612\"],\n  [\"62\u00d759=3658\", \"24\u00d716=384\"],\n  [\"12\u00d733=396\", \"54\u00d763=3402\"],\n  [\"72\u00d738=2736\", \"11\u00d724=264\"],\n  [\"51\u00d769=3519\", \"91\u00d738=3458\"],\n  [\"41\u00d774=3034\", \"39\u00d729=1131\"],\n  [\"89\u00d715=1335\", \"77\u00d782=6314\"],\n  [\"22\u00d764=1408\", \"74\u00d738=2812\"],\n  [\"56\u00d789=4984\", \"29\u00d730=870\"],\n  [\"60\u00d756=3360\", \"89\u00d723=2047\"],\n  [\"47\u00d757=2679\", \"85\u00d783=7055\"],\n  [\"29\u00d798=2842\", \"34\u00d758=1972\"],\n  [\"61\u00d735=2135\", \"93\u00d730=2790\"],\n  [\"24\u00d715=360\", \"88\u00d728=2464\"],\n  [\"42\u00d758=2436\", \"62\u00d757=3534\"],\n  [\"49\u00d725=1225\", \"47\u00d746=2162\"],\n  [\"87\u00d730=2610\", \"15\u00d733=495\"],\n  [\"62\u00d779=4898\", \"29\u00d757=1653\"],\n  [\"28\u00d744=1232\", \"81\u00d788=7128\"],\n  [\"34\u00d794=3196\", \"89\u00d756=4984\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    console.log(\"WARNING: no match found for '\" + oldText + \"'\");\n    continue;\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each two-digit multiplication problem's text\n# with the updated values, matching the authoritative diff 1:1.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @('2024-12-15 Sunday', '2024-12-16 Monday'),\n    @('56\u00d794=5264', '39\u00d786=3354'),\n    @('32\u00d762=1984', '77\u00d775=5775'),\n    @('46\u00d746=2116', '55\u00d720=1100'),\n    @('34\u00d730=1020', '12\u00d711=132'),\n    @('25\u00d744=1100', '61\u00d724=1464'),\n    @('31\u00d732=992', '31\u00d752=1612'),\n    @('62\u00d759=3658', '24\u00d716=384'),\n    @('12\u00d733=396', '54\u00d763=3402'),\n    @('72\u00d738=2736', '11\u00d724=264'),\n    @('51\u00d769=3519', '91\u00d738=3458'),\n    @('41\u00d774=3034', '39\u00d729=1131'),\n    @('89\u00d715=1335', '77\u00d782=6314'),\n    @('22\u00d764=1408', '74\u00d738=2812'),\n    @('56\u00d789=4984', '29\u00d730=870'),\n    @('60\u00d756=3360', '89\u00d723=2047'),\n    @('47\u00d757=2679', '85\u00d783=7055'),\n    @('29\u00d798=2842', '34\u00d758=1972'),\n    @('61\u00d735=2135', '93\u00d730=2790'),\n    @('24\u00d715=360', '88\u00d728=2464'),\n    @('42\u00d758=2436', '62\u00d757=3534'),\n    @('49\u00d725=1225', '47\u00d746=2162'),\n    @('87\u00d730=2610', '15\u00d733=495'),\n    @('62\u00d779=4898', '29\u00d757=1653'),\n    @('28\u00d744=1232', '81\u00d788=7128'),\n    @('34\u00d794=3196', '89\u00d756=4984'),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $found) {\n        Write-Output \"WARNING: no match found for '$oldText'\"\n    }\n}\n"}
